# "Added names to user stories"
#
# Slide 3 ("User Stories 1"):
#   - Title text:  "User Stories 1 "  -> "User Stories - Jared "
#   - Body placeholder moved down (Top: 1430100 EMU -> 1582500 EMU)
#   - The second bullet ("...tracks that move... (In Progress)(4)") is removed
#     entirely (its paragraph is deleted / merged away).
# Slide 4 ("User Stories 2"):
#   - Title text:  "User Stories 2" -> "User Stories - Jennifer"
# Slide 5 ("User Stories 3"):
#   - Title text:  "User Stories 3" -> "User Stories - Christian"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 : "User Stories 1" -> "User Stories - Jared"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Title placeholder is shape 1 on this slide.
$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "User Stories - Jared "

# Body placeholder is shape 2 on this slide.
$body3 = $s3.Shapes.Item(2)

# Move the body box down: new Top = 1582500 EMU (EMU / 12700 = points).
$body3.Top = 1582500 / 12700

# Remove the second bullet paragraph in its entirety:
#   "As a user of the system, I want to be able to manually add tracks that
#    move by clicking a button on the map, so that I can make predictions
#    about a track's course of collision. (In Progress)(4)"
# That paragraph (plus its trailing paragraph mark) starts right after the
# first paragraph's mark (character 194) and is 194 characters long
# (193 characters of run text + its own trailing paragraph mark), so
# deleting Characters(195, 194) removes the run text AND merges away the
# paragraph break, joining paragraph 1 directly to what was paragraph 3.
$tr3 = $body3.TextFrame.TextRange
$deadPara = $tr3.Characters(195, 194)
$deadPara.Delete()

# ---------------------------------------------------------------------------
# Slide 4 : "User Stories 2" -> "User Stories - Jennifer"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$title4 = $s4.Shapes.Item(1)
$title4.TextFrame.TextRange.Text = "User Stories - Jennifer"

# ---------------------------------------------------------------------------
# Slide 5 : "User Stories 3" -> "User Stories - Christian"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
# On this slide the title placeholder happens to be shape 2 (body is shape 1).
$title5 = $s5.Shapes.Item(2)
$title5.TextFrame.TextRange.Text = "User Stories - Christian"
